# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) counts on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 16976
$ws1.Range("F9").Value  = 12
$ws1.Range("F13").Value = 11746
$ws1.Range("F15").Value = 16
$ws1.Range("F16").Value = 1438
$ws1.Range("F17").Value = 4672
$ws1.Range("F18").Value = 480
$ws1.Range("F19").Value = 17
$ws1.Range("F21").Value = 73
$ws1.Range("F25").Value = 32

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 16976
$ws4.Range("F10").Value = 12
$ws4.Range("F16").Value = 11746
$ws4.Range("F18").Value = 16
$ws4.Range("F19").Value = 1438
$ws4.Range("F20").Value = 4672
$ws4.Range("F21").Value = 480
$ws4.Range("F22").Value = 17
$ws4.Range("F24").Value = 73
$ws4.Range("F28").Value = 32
